$d = $word.ActiveDocument

# --- Step 1: paragraph "Библиотека для разработки игр Pygame." ---
# its last run (just the trailing period, its own run with distinct
# formatting) changes text from "." to ";" while keeping the run's
# original rsid metadata intact.
$pPygame = $d.Paragraphs(14)
$pygameRange = $pPygame.Range
$periodEnd = $pygameRange.End
$periodStart = $periodEnd - 2
$periodEndExclParaMark = $periodEnd - 1
$periodRange = $d.Range($periodStart, $periodEndExclParaMark)

$xmlPkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlPkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$semiXml = $xmlPkgHeader
$semiXml += '<w:p><w:r w:rsidRPr="00254266"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr><w:t>;</w:t></w:r></w:p>'
$semiXml += $xmlPkgFooter
$periodRange.InsertXML($semiXml)

# --- Step 2: insert three new bulleted paragraphs after it, replacing
# the blank paragraph that used to follow, with: "Библиотека sys;",
# "Библиотека math;", "Библиотека os." (sys/math/os bold + en-US). ---
$pPygame = $d.Paragraphs(14)
$pPygame.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(15)
$newRange = $newPara.Range

$libXml = $xmlPkgHeader

$libXml += '<w:p>'
$libXml += '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr></w:pPr>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Библиотека </w:t></w:r>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>sys</w:t></w:r>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>;</w:t></w:r>'
$libXml += '</w:p>'

$libXml += '<w:p>'
$libXml += '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr></w:pPr>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Библиотека </w:t></w:r>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>math</w:t></w:r>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>;</w:t></w:r>'
$libXml += '</w:p>'

$libXml += '<w:p>'
$libXml += '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr></w:pPr>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Библиотека </w:t></w:r>'
$libXml += '<w:proofErr w:type="spellStart"/>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>os</w:t></w:r>'
$libXml += '<w:proofErr w:type="spellEnd"/>'
$libXml += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>'
$libXml += '</w:p>'

$libXml += $xmlPkgFooter

$newRange.InsertXML($libXml)

Write-Output "ok"
